$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 data: iteration 6, differences "56/200", training data 1800
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "56/200"
$ws.Cells.Item(7, 3).Value = 1800

# Update selection to reflect the new active cell state (C8)
$ws.Range("C8").Select()
